# Update column C (Förändrad) date from 2023-09-09 (45178) to 2023-09-10 (45179)
# for all data rows (2 through 57) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 57
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45178) {
        $cell.Value2 = 45179
    }
}
